$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 85  # H12
$ws.Cells.Item(12, 9).Value = 85  # I12
$ws.Cells.Item(12, 10).Value = 85  # J12
$ws.Cells.Item(12, 11).Value = 85  # K12
$ws.Cells.Item(12, 12).Value = 85  # L12
$ws.Cells.Item(12, 13).Value = 85  # M12
$ws.Cells.Item(12, 14).Value = -425  # N12
$ws.Cells.Item(21, 8).Value = 42479.5  # H21
$ws.Cells.Item(21, 9).Value = 54972.668  # I21
$ws.Cells.Item(21, 10).Value = 5000  # J21
$ws.Cells.Item(21, 11).Value = 54972.668  # K21
$ws.Cells.Item(21, 12).Value = 5000  # L21
$ws.Cells.Item(21, 13).Value = -54504.668  # M21
$ws.Cells.Item(21, 14).Value = -5936  # N21
$ws.Cells.Item(23, 8).Value = 42479.5  # H23
$ws.Cells.Item(23, 9).Value = 54972.668  # I23
$ws.Cells.Item(23, 10).Value = 5000  # J23
$ws.Cells.Item(23, 11).Value = 54972.668  # K23
$ws.Cells.Item(23, 12).Value = 5000  # L23
$ws.Cells.Item(23, 13).Value = -54738.668  # M23
$ws.Cells.Item(23, 14).Value = -5468  # N23
$ws.Cells.Item(58, 8).Value = 2692.4211  # H58
$ws.Cells.Item(58, 9).Value = 262.7143  # I58
$ws.Cells.Item(58, 10).Value = 4109.75  # J58
$ws.Cells.Item(58, 11).Value = 788.1428999999999  # K58
$ws.Cells.Item(58, 12).Value = 12329.25  # L58
$ws.Cells.Item(58, 13).Value = -638.1428999999999  # M58
$ws.Cells.Item(58, 14).Value = -12629.25  # N58
$ws.Cells.Item(61, 8).Value = 2231.818  # H61
$ws.Cells.Item(61, 9).Value = 154.125  # I61
$ws.Cells.Item(61, 10).Value = 7772.3335  # J61
$ws.Cells.Item(61, 11).Value = 462.375  # K61
$ws.Cells.Item(61, 12).Value = 23317.0005  # L61
$ws.Cells.Item(61, 13).Value = -290.375  # M61
$ws.Cells.Item(61, 14).Value = -23661.0005  # N61
$ws.Cells.Item(69, 8).Value = 5603.75  # H69
$ws.Cells.Item(69, 9).Value = 6000  # I69
$ws.Cells.Item(69, 10).Value = 5366  # J69
$ws.Cells.Item(69, 11).Value = 18000  # K69
$ws.Cells.Item(69, 12).Value = 16098  # L69
$ws.Cells.Item(69, 13).Value = -17126  # M69
$ws.Cells.Item(69, 14).Value = -17846  # N69
$ws.Cells.Item(70, 8).Value = 948.3871  # H70
$ws.Cells.Item(70, 9).Value = 954.7917  # I70
$ws.Cells.Item(70, 10).Value = 926.4286  # J70
$ws.Cells.Item(70, 11).Value = 2864.3751  # K70
$ws.Cells.Item(70, 12).Value = 2779.2858  # L70
$ws.Cells.Item(70, 13).Value = -2594.3751  # M70
$ws.Cells.Item(70, 14).Value = -3319.2858  # N70
$ws.Cells.Item(72, 8).Value = 5603.75  # H72
$ws.Cells.Item(72, 9).Value = 6000  # I72
$ws.Cells.Item(72, 10).Value = 5366  # J72
$ws.Cells.Item(72, 11).Value = 54000  # K72
$ws.Cells.Item(72, 12).Value = 48294  # L72
$ws.Cells.Item(72, 13).Value = -49632  # M72
$ws.Cells.Item(72, 14).Value = -57030  # N72
$ws.Cells.Item(73, 8).Value = 948.3871  # H73
$ws.Cells.Item(73, 9).Value = 954.7917  # I73
$ws.Cells.Item(73, 10).Value = 926.4286  # J73
$ws.Cells.Item(73, 11).Value = 2864.3751  # K73
$ws.Cells.Item(73, 12).Value = 2779.2858  # L73
$ws.Cells.Item(73, 13).Value = -1928.3751  # M73
$ws.Cells.Item(73, 14).Value = -4651.2858  # N73
$ws.Cells.Item(80, 8).Value = 610.1111  # H80
$ws.Cells.Item(80, 9).Value = 1186.6666  # I80
$ws.Cells.Item(80, 10).Value = 494.8  # J80
$ws.Cells.Item(80, 11).Value = 3559.9998  # K80
$ws.Cells.Item(80, 12).Value = 1484.4  # L80
$ws.Cells.Item(80, 13).Value = -2561.9998  # M80
$ws.Cells.Item(80, 14).Value = -3480.4  # N80
$ws.Cells.Item(82, 8).Value = 3387.875  # H82
$ws.Cells.Item(82, 9).Value = 400.54544  # I82
$ws.Cells.Item(82, 10).Value = 9960  # J82
$ws.Cells.Item(82, 11).Value = 1201.63632  # K82
$ws.Cells.Item(82, 12).Value = 29880  # L82
$ws.Cells.Item(82, 13).Value = -795.6363200000001  # M82
$ws.Cells.Item(82, 14).Value = -30692  # N82
$ws.Cells.Item(83, 8).Value = 610.1111  # H83
$ws.Cells.Item(83, 9).Value = 1186.6666  # I83
$ws.Cells.Item(83, 10).Value = 494.8  # J83
$ws.Cells.Item(83, 11).Value = 10679.9994  # K83
$ws.Cells.Item(83, 12).Value = 4453.2  # L83
$ws.Cells.Item(83, 13).Value = -5687.999400000001  # M83
$ws.Cells.Item(83, 14).Value = -14437.2  # N83
$ws.Cells.Item(85, 8).Value = 3387.875  # H85
$ws.Cells.Item(85, 9).Value = 400.54544  # I85
$ws.Cells.Item(85, 10).Value = 9960  # J85
$ws.Cells.Item(85, 11).Value = 1201.63632  # K85
$ws.Cells.Item(85, 12).Value = 29880  # L85
$ws.Cells.Item(85, 13).Value = 202.3636799999999  # M85
$ws.Cells.Item(85, 14).Value = -32688  # N85
$ws.Cells.Item(100, 8).Value = 1591.4117  # H100
$ws.Cells.Item(100, 9).Value = 1546.7142  # I100
$ws.Cells.Item(100, 10).Value = 1800  # J100
$ws.Cells.Item(100, 11).Value = 1546.7142  # K100
$ws.Cells.Item(100, 12).Value = 1800  # L100
$ws.Cells.Item(100, 13).Value = -1005.7142  # M100
$ws.Cells.Item(100, 14).Value = -2882  # N100
$ws.Cells.Item(129, 8).Value = 1201.9531  # H129
$ws.Cells.Item(129, 9).Value = 906.0909  # I129
$ws.Cells.Item(129, 10).Value = 1263.3585  # J129
$ws.Cells.Item(129, 11).Value = 2718.2727  # K129
$ws.Cells.Item(129, 12).Value = 3790.0755  # L129
$ws.Cells.Item(129, 13).Value = 2281.7273  # M129
$ws.Cells.Item(129, 14).Value = -13790.0755  # N129
$ws.Cells.Item(137, 8).Value = 2560.4211  # H137
$ws.Cells.Item(137, 9).Value = 4632.25  # I137
$ws.Cells.Item(137, 10).Value = 1053.6364  # J137
$ws.Cells.Item(137, 11).Value = 13896.75  # K137
$ws.Cells.Item(137, 12).Value = 3160.9092  # L137
$ws.Cells.Item(137, 13).Value = -11346.75  # M137
$ws.Cells.Item(137, 14).Value = -8260.9092  # N137
$ws.Cells.Item(138, 8).Value = 2098.875  # H138
$ws.Cells.Item(138, 9).Value = 1725.6  # I138
$ws.Cells.Item(138, 10).Value = 2338.1538  # J138
$ws.Cells.Item(138, 11).Value = 5176.799999999999  # K138
$ws.Cells.Item(138, 12).Value = 7014.4614  # L138
$ws.Cells.Item(138, 13).Value = -36.79999999999927  # M138
$ws.Cells.Item(138, 14).Value = -17294.4614  # N138

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 8810.111999999999  # H32
$ws.Cells.Item(32, 9).Value = 7676.1626  # I32
$ws.Cells.Item(32, 10).Value = 16936.75  # J32
$ws.Cells.Item(32, 11).Value = 7676.1626  # K32
$ws.Cells.Item(32, 12).Value = 16936.75  # L32
$ws.Cells.Item(32, 13).Value = -7389.1626  # M32
$ws.Cells.Item(32, 14).Value = -17510.75  # N32
$ws.Cells.Item(96, 8).Value = 16836  # H96
$ws.Cells.Item(96, 10).Value = 16836  # J96
$ws.Cells.Item(96, 12).Value = 16836  # L96
$ws.Cells.Item(96, 14).Value = -22328  # N96
$ws.Cells.Item(102, 8).Value = 2278.4211  # H102
$ws.Cells.Item(102, 9).Value = 1969.375  # I102
$ws.Cells.Item(102, 10).Value = 3926.6667  # J102
$ws.Cells.Item(102, 11).Value = 1969.375  # K102
$ws.Cells.Item(102, 12).Value = 3926.6667  # L102
$ws.Cells.Item(102, 13).Value = -347.375  # M102
$ws.Cells.Item(102, 14).Value = -7170.6667  # N102

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(43, 8).Value = 9392.799999999999  # H43
$ws.Cells.Item(43, 10).Value = 9392.799999999999  # J43
$ws.Cells.Item(43, 12).Value = 9392.799999999999  # L43
$ws.Cells.Item(43, 14).Value = -9760.799999999999  # N43
$ws.Cells.Item(101, 8).Value = 9392.799999999999  # H101
$ws.Cells.Item(101, 10).Value = 9392.799999999999  # J101
$ws.Cells.Item(101, 12).Value = 9392.799999999999  # L101
$ws.Cells.Item(101, 14).Value = -15882.8  # N101

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(141, 8).Value = 3670  # H141
$ws.Cells.Item(141, 9).Value = 2032.5  # I141
$ws.Cells.Item(141, 10).Value = 4980  # J141
$ws.Cells.Item(141, 11).Value = 6097.5  # K141
$ws.Cells.Item(141, 12).Value = 14940  # L141
$ws.Cells.Item(141, 13).Value = -917.5  # M141
$ws.Cells.Item(141, 14).Value = -25300  # N141

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(13, 8).Value = 576.25  # H13
$ws.Cells.Item(13, 9).Value = 726.25  # I13
$ws.Cells.Item(13, 10).Value = 426.25  # J13
$ws.Cells.Item(13, 11).Value = 726.25  # K13
$ws.Cells.Item(13, 12).Value = 426.25  # L13
$ws.Cells.Item(13, 13).Value = -587.25  # M13
$ws.Cells.Item(13, 14).Value = -704.25  # N13
$ws.Cells.Item(41, 8).Value = 9980  # H41
$ws.Cells.Item(41, 9).Value = 0  # I41
$ws.Cells.Item(41, 10).Value = 9980  # J41
$ws.Cells.Item(41, 11).Value = 0  # K41
$ws.Cells.Item(41, 12).Value = 9980  # L41
$ws.Cells.Item(41, 13).ClearContents()  # M41
$ws.Cells.Item(41, 14).Value = -10690  # N41
$ws.Cells.Item(99, 8).Value = 10785.8125  # H99
$ws.Cells.Item(99, 9).Value = 7535.727  # I99
$ws.Cells.Item(99, 11).Value = 7535.727  # K99
$ws.Cells.Item(99, 13).Value = -5289.727  # M99
$ws.Cells.Item(102, 8).Value = 9245.171  # H102
$ws.Cells.Item(102, 9).Value = 8929.741  # I102
$ws.Cells.Item(102, 10).Value = 9853.5  # J102
$ws.Cells.Item(102, 11).Value = 8929.741  # K102
$ws.Cells.Item(102, 12).Value = 9853.5  # L102
$ws.Cells.Item(102, 13).Value = -7307.741  # M102
$ws.Cells.Item(102, 14).Value = -13097.5  # N102
$ws.Cells.Item(132, 8).Value = 2087619.4  # H132
$ws.Cells.Item(132, 9).Value = 4632942.5  # I132
$ws.Cells.Item(132, 10).Value = 5082.364  # J132
$ws.Cells.Item(132, 11).Value = 13898827.5  # K132
$ws.Cells.Item(132, 12).Value = 15247.092  # L132
$ws.Cells.Item(132, 13).Value = -13896297.5  # M132
$ws.Cells.Item(132, 14).Value = -20307.092  # N132

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(106, 8).Value = 33000  # H106
$ws.Cells.Item(106, 10).Value = 33000  # J106
$ws.Cells.Item(106, 12).Value = 33000  # L106
$ws.Cells.Item(106, 14).Value = -35524  # N106

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 1633.1538  # H81
$ws.Cells.Item(81, 10).Value = 1976.6666  # J81
$ws.Cells.Item(81, 12).Value = 3953.3332  # L81
$ws.Cells.Item(81, 14).Value = -6075.3332  # N81
$ws.Cells.Item(84, 8).Value = 1633.1538  # H84
$ws.Cells.Item(84, 10).Value = 1976.6666  # J84
$ws.Cells.Item(84, 12).Value = 19766.666  # L84
$ws.Cells.Item(84, 14).Value = -30374.666  # N84
$ws.Cells.Item(94, 8).Value = 30330  # H94
$ws.Cells.Item(94, 10).Value = 30330  # J94
$ws.Cells.Item(94, 12).Value = 30330  # L94
$ws.Cells.Item(94, 14).Value = -32132  # N94
$ws.Cells.Item(101, 8).Value = 8701  # H101
$ws.Cells.Item(101, 10).Value = 8701  # J101
$ws.Cells.Item(101, 12).Value = 8701  # L101
$ws.Cells.Item(101, 14).Value = -15191  # N101
